$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header column F1 from "ukuran" to "warna"
$ws.Range("F1").Value = "warna"

# Update the selected cell/active cell to H6
$ws.Range("H6").Select()
